# "modified the format to fit the event table"
#
# The sheet's header row is reshuffled: a new lowercase "summary" header
# becomes column A, and the previous headers (dtstart/dtend/private/userid)
# shift one column to the right (B:E). The "private" column's values switch
# from the literal text "t" to real boolean TRUE/FALSE. The active selection
# moves from E3 to D4, and the saved window position changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: insert "summary" in A1 and push dtstart/dtend/private/userid
# one column over.
$ws.Range("A1").Value = "summary"
$ws.Range("B1").Value = "dtstart"
$ws.Range("C1").Value = "dtend"
$ws.Range("D1").Value = "private"
$ws.Range("E1").Value = "userid"

# "private" column becomes real booleans instead of the text "t".
$ws.Range("D2").Value = $true
$ws.Range("D3").Value = $false

# Active cell/selection moves from E3 to D4.
$ws.Range("D4").Select()

# Saved window position (best effort - host may not persist this).
$excel.ActiveWindow.Left = 1420
$excel.ActiveWindow.Top = 1740
